$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2293
$ws.Range("J32").Value = 1964.7273
$ws.Range("L32").Value = 1964.7273
$ws.Range("N32").Value = -2616.7273

$ws.Range("H40").Value = 1426.3334
$ws.Range("I40").Value = 1251.4814
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 1251.4814
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1076.4814
$ws.Range("N40").Value = -3350

$ws.Range("H86").Value = 7731.8335
$ws.Range("I86").Value = 9165.666999999999
$ws.Range("J86").Value = 6298
$ws.Range("K86").Value = 9165.666999999999
$ws.Range("L86").Value = 6298
$ws.Range("M86").Value = -8042.666999999999
$ws.Range("N86").Value = -8544

$ws.Range("H88").Value = 1580.1111
$ws.Range("J88").Value = 1051.5
$ws.Range("L88").Value = 1051.5
$ws.Range("N88").Value = -1863.5

$ws.Range("H89").Value = 7731.8335
$ws.Range("I89").Value = 9165.666999999999
$ws.Range("J89").Value = 6298
$ws.Range("K89").Value = 45828.335
$ws.Range("L89").Value = 31490
$ws.Range("M89").Value = -40212.335
$ws.Range("N89").Value = -42722

$ws.Range("H91").Value = 1580.1111
$ws.Range("J91").Value = 1051.5
$ws.Range("L91").Value = 1051.5
$ws.Range("N91").Value = -3859.5

$ws.Range("H138").Value = 4878.6772
$ws.Range("I138").Value = 11111
$ws.Range("J138").Value = 4670.933
$ws.Range("K138").Value = 33333
$ws.Range("L138").Value = 14012.799
$ws.Range("M138").Value = -28193
$ws.Range("N138").Value = -24292.799

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 141081.81
$ws.Range("J32").Value = 17947.861
$ws.Range("L32").Value = 17947.861
$ws.Range("N32").Value = -18521.861

$ws.Range("H61").Value = 1518751
$ws.Range("I61").Value = 3605.9038
$ws.Range("K61").Value = 3605.9038
$ws.Range("M61").Value = -3393.9038

$ws.Range("H74").Value = 1993115.4
$ws.Range("I74").Value = 2650335.2
$ws.Range("K74").Value = 2650335.2
$ws.Range("M74").Value = -2649461.2

$ws.Range("H77").Value = 1993115.4
$ws.Range("I77").Value = 2650335.2
$ws.Range("K77").Value = 13251676
$ws.Range("M77").Value = -13247308

$ws.Range("H122").Value = 1538.375
$ws.Range("I122").Value = 1308.7693
$ws.Range("K122").Value = 3926.3079
$ws.Range("M122").Value = -1476.3079

$ws.Range("H136").Value = 1518751
$ws.Range("I136").Value = 3605.9038
$ws.Range("K136").Value = 10817.7114
$ws.Range("M136").Value = -8267.7114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1002.2
$ws.Range("I107").Value = 1002.2
$ws.Range("K107").Value = 1002.2
$ws.Range("M107").Value = 917.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1685.1666
$ws.Range("I25").Value = 1952.875
$ws.Range("J25").Value = 1149.75
$ws.Range("K25").Value = 1952.875
$ws.Range("L25").Value = 1149.75
$ws.Range("M25").Value = -1778.875
$ws.Range("N25").Value = -1497.75

$ws.Range("H31").Value = 2648710.8
$ws.Range("I31").Value = 3089273.8
$ws.Range("J31").Value = 5333.3335
$ws.Range("K31").Value = 3089273.8
$ws.Range("L31").Value = 5333.3335
$ws.Range("M31").Value = -3088978.8
$ws.Range("N31").Value = -5923.3335

$ws.Range("H32").Value = 2124.75
$ws.Range("I32").Value = 2333
$ws.Range("K32").Value = 2333
$ws.Range("M32").Value = -2017

$ws.Range("H33").Value = 3323.3333
$ws.Range("I33").Value = 3323.3333
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 3323.3333
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -2944.3333
$ws.Range("N33").Value = $null

$ws.Range("H34").Value = 2648710.8
$ws.Range("I34").Value = 3089273.8
$ws.Range("J34").Value = 5333.3335
$ws.Range("K34").Value = 3089273.8
$ws.Range("L34").Value = 5333.3335
$ws.Range("M34").Value = -3089071.8
$ws.Range("N34").Value = -5737.3335

$ws.Range("H58").Value = 4913894.5
$ws.Range("I58").Value = 8805.5
$ws.Range("J58").Value = 6423153
$ws.Range("K58").Value = 8805.5
$ws.Range("L58").Value = 6423153
$ws.Range("M58").Value = -8602.5
$ws.Range("N58").Value = -6423559

$ws.Range("H105").Value = 6147.909
$ws.Range("J105").Value = 2264.5
$ws.Range("L105").Value = 2264.5
$ws.Range("N105").Value = -5758.5

$ws.Range("H132").Value = 1396.7142
$ws.Range("I132").Value = 1396.7142
$ws.Range("K132").Value = 4190.142599999999
$ws.Range("M132").Value = -1660.142599999999

$ws.Range("H134").Value = 1704.0526
$ws.Range("I134").Value = 1576.5
$ws.Range("K134").Value = 4729.5
$ws.Range("M134").Value = -2194.5

$ws.Range("H136").Value = 4913894.5
$ws.Range("I136").Value = 8805.5
$ws.Range("J136").Value = 6423153
$ws.Range("K136").Value = 26416.5
$ws.Range("L136").Value = 19269459
$ws.Range("M136").Value = -23866.5
$ws.Range("N136").Value = -19274559

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2565530.5
$ws.Range("J5").Value = 2328559.5
$ws.Range("L5").Value = 6985678.5
$ws.Range("N5").Value = -6985902.5

$ws.Range("H23").Value = 84.42104999999999
$ws.Range("J23").Value = 86.111115
$ws.Range("L23").Value = 258.333345
$ws.Range("N23").Value = -728.333345

$ws.Range("H133").Value = 16166.333
$ws.Range("I133").Value = 16166.333
$ws.Range("K133").Value = 48498.999
$ws.Range("M133").Value = -43438.999

$ws.Range("H135").Value = 2565530.5
$ws.Range("J135").Value = 2328559.5
$ws.Range("L135").Value = 20957035.5
$ws.Range("N135").Value = -20962105.5

$ws.Range("H139").Value = 5482.909
$ws.Range("I139").Value = 3152.4783
$ws.Range("K139").Value = 9457.4349
$ws.Range("M139").Value = -4317.4349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 13999.889
$ws.Range("I18").Value = 13999.889
$ws.Range("K18").Value = 13999.889
$ws.Range("M18").Value = -13706.889

$ws.Range("H80").Value = 2583.8
$ws.Range("I80").Value = 2604.75
$ws.Range("K80").Value = 2604.75
$ws.Range("M80").Value = -1606.75

$ws.Range("H83").Value = 2583.8
$ws.Range("I83").Value = 2604.75
$ws.Range("K83").Value = 13023.75
$ws.Range("M83").Value = -8031.75

$ws.Range("H126").Value = 7348.5713
$ws.Range("I126").Value = 8498.182000000001
$ws.Range("K126").Value = 25494.546
$ws.Range("M126").Value = -23024.546

$ws.Range("H132").Value = 15977.65
$ws.Range("I132").Value = 13650.467
$ws.Range("K132").Value = 40951.401
$ws.Range("M132").Value = -38421.401

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1112.1111
$ws.Range("J22").Value = 1214.9286
$ws.Range("L22").Value = 1214.9286
$ws.Range("N22").Value = -1804.9286

$ws.Range("H27").Value = 1112.1111
$ws.Range("J27").Value = 1214.9286
$ws.Range("L27").Value = 1214.9286
$ws.Range("N27").Value = -1428.9286

$ws.Range("H82").Value = 2818.05
$ws.Range("J82").Value = 2304.5
$ws.Range("L82").Value = 2304.5
$ws.Range("N82").Value = -3026.5

$ws.Range("H85").Value = 2818.05
$ws.Range("J85").Value = 2304.5
$ws.Range("L85").Value = 2304.5
$ws.Range("N85").Value = -4800.5

$ws.Range("H122").Value = 3298.7778
$ws.Range("I122").Value = 2529.1428
$ws.Range("K122").Value = 7587.428400000001
$ws.Range("M122").Value = -5137.428400000001

$ws.Range("H132").Value = 5317452.5
$ws.Range("I132").Value = 8990641
$ws.Range("K132").Value = 26971923
$ws.Range("M132").Value = -26969393

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 30000
$ws.Range("I32").Value = 30000
$ws.Range("K32").Value = 30000
$ws.Range("M32").Value = -29683

$ws.Range("H62").Value = 8124.3335
$ws.Range("I62").Value = 6483.5
$ws.Range("J62").Value = 8721
$ws.Range("K62").Value = 6483.5
$ws.Range("L62").Value = 8721
$ws.Range("M62").Value = -5859.5
$ws.Range("N62").Value = -9969

$ws.Range("H65").Value = 8124.3335
$ws.Range("I65").Value = 6483.5
$ws.Range("J65").Value = 8721
$ws.Range("K65").Value = 32417.5
$ws.Range("L65").Value = 43605
$ws.Range("M65").Value = -29297.5
$ws.Range("N65").Value = -49845

$ws.Range("H113").Value = 899.625
$ws.Range("I113").Value = 944.1111
$ws.Range("K113").Value = 2832.3333
$ws.Range("M113").Value = -662.3332999999998

$ws.Range("H136").Value = 5939501
$ws.Range("I136").Value = 1359557.5
$ws.Range("J136").Value = 22223744
$ws.Range("K136").Value = 4078672.5
$ws.Range("L136").Value = 66671232
$ws.Range("M136").Value = -4076122.5
$ws.Range("N136").Value = -66676332
